$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting existing rows 97:151 down to 98:152
$ws.Rows("97:97").Insert()

# Populate the newly inserted row 97 with the new record's data
$ws.Range("A97").Value = 4
$ws.Range("B97").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C97").Value = "Los Lagos"
$ws.Range("D97").Value = 45097
$ws.Range("D97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E97").Value = 10
$ws.Range("F97").Value = 100112031
$ws.Range("G97").Value = "Poroto verde"
$ws.Range("H97").Value = "Magnum"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 40
$ws.Range("K97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("M97").Value = 30000
$ws.Range("N97").Value = "`$/malla 25 kilos"
$ws.Range("O97").Value = "Perú"
$ws.Range("P97").Value = 1200
$ws.Range("Q97").Value = 25
$ws.Range("R97").Value = "Hortaliza"
